$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H2" = 519.3333
    "I2" = 96.42856999999999
    "K2" = 96.42856999999999
    "M2" = 16.57143000000001
    "H33" = 619.64703
    "I33" = 288.92856
    "J33" = 2163
    "K33" = 288.92856
    "L33" = 2163
    "M33" = -59.92856
    "N33" = -2621
    "H98" = 1867.7
    "I98" = 742.1111
    "K98" = 742.1111
    "M98" = 755.8889
    "H106" = 6902.5293
    "I106" = 6902.5293
    "K106" = 6902.5293
    "M106" = -6271.5293
    "H122" = 1867.7
    "I122" = 742.1111
    "K122" = 2226.3333
    "M122" = 223.6667000000002
    "H137" = 2220.8235
    "I137" = 2234.625
    "K137" = 6703.875
    "M137" = -4153.875
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H2" = 1802.5294
    "I2" = 1493.7742
    "K2" = 1493.7742
    "M2" = -1380.7742
    "H31" = 13694
    "I31" = 13694
    "K31" = 13694
    "M31" = -13400
    "H61" = 4821.4136
    "I61" = 3744.0454
    "K61" = 3744.0454
    "M61" = -3532.0454
    "H74" = 2206.3
    "I74" = 2279.12
    "J74" = 1842.2
    "K74" = 2279.12
    "L74" = 1842.2
    "M74" = -1405.12
    "N74" = -3590.2
    "H77" = 2206.3
    "I77" = 2279.12
    "J77" = 1842.2
    "K77" = 11395.6
    "L77" = 9211
    "M77" = -7027.599999999999
    "N77" = -17947
    "H116" = 1802.5294
    "I116" = 1493.7742
    "K116" = 1493.7742
    "M116" = 800.2257999999999
    "H122" = 3095.4167
    "I122" = 2480.9805
    "K122" = 7442.941500000001
    "M122" = -4992.941500000001
    "H132" = 1881.75
    "I132" = 1545.5151
    "K132" = 4636.5453
    "M132" = -2106.5453
    "H136" = 4821.4136
    "I136" = 3744.0454
    "K136" = 11232.1362
    "M136" = -8682.136200000001
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H3" = 1802.5294
    "I3" = 1493.7742
    "K3" = 1493.7742
    "M3" = -1379.7742
    "H86" = 2485.125
    "I86" = 2418.2222
    "K86" = 2418.2222
    "M86" = -1295.2222
    "H89" = 2485.125
    "I89" = 2418.2222
    "K89" = 12365
    "M89" = -6475.111000000001
    "H134" = 4122.2
    "J134" = 5991.5
    "L134" = 17974.5
    "N134" = -23044.5
    "H61" = 699.5
    "I61" = 699.5
    "K61" = 699.5
    "M61" = -497.5
    "H113" = 699.5
    "I113" = 699.5
    "K113" = 699.5
    "M113" = 1470.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H16" = 10114.4375
    "I16" = 712.1818
    "J16" = 30799.4
    "K16" = 712.1818
    "L16" = 30799.4
    "M16" = -425.1818
    "N16" = -31373.4
    "H20" = 94250
    "J20" = 94250
    "L20" = 94250
    "N20" = -94722
    "H30" = 94250
    "J30" = 94250
    "L30" = 94250
    "N30" = -94432
    "H62" = 6690.5625
    "J62" = 6211
    "L62" = 6211
    "N62" = -7459
    "H65" = 6690.5625
    "J65" = 6211
    "L65" = 31055
    "N65" = -37295
    "H105" = 6599.5557
    "I105" = 4999.5713
    "K105" = 4999.5713
    "M105" = -3252.5713
    "H113" = 10114.4375
    "I113" = 712.1818
    "J113" = 30799.4
    "K113" = 712.1818
    "L113" = 30799.4
    "M113" = 1457.8182
    "N113" = -35139.4
    "H122" = 3147.5862
    "I122" = 3291.762
    "J122" = 2769.125
    "K122" = 9875.286
    "L122" = 8307.375
    "M122" = -7425.286
    "N122" = -13207.375
    "H128" = 94250
    "J128" = 94250
    "L128" = 94250
    "N128" = -104210
    "H141" = 294247.5
    "J141" = 294247.5
    "L141" = 294247.5
    "N141" = -304607.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H4" = 250422600
    "I4" = 333563460
    "J4" = 1000000
    "K4" = 1000690380
    "L4" = 3000000
    "M4" = -1000690268
    "N4" = -3000224
    "H45" = 766.3333
    "J45" = 766.3333
    "L45" = 2298.9999
    "N45" = -3362.9999
    "H108" = 1853.3334
    "I108" = 1824
    "K108" = 5472
    "M108" = -2592
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H97" = 436.14285
    "I97" = 407.95
    "K97" = 407.95
    "M97" = 88.05000000000001
    "H136" = 29607.334
    "J136" = 29607.334
    "L136" = 88822.00199999999
    "N136" = -93922.00199999999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H46" = 5100
    "I46" = 0
    "K46" = 0
    "H68" = 4496.75
    "I68" = 3997.5
    "K68" = 3997.5
    "M68" = -3248.5
    "H71" = 4496.75
    "I71" = 3997.5
    "K71" = 19987.5
    "M71" = -16243.5
    "H93" = 1559.9412
    "I93" = 1482.2222
    "J93" = 1647.375
    "K93" = 1482.2222
    "L93" = 1647.375
    "M93" = -234.2221999999999
    "N93" = -4143.375
    "H100" = 1998
    "I100" = 1998
    "K100" = 1998
    "M100" = -1457
    "H102" = 98780.5
    "J102" = 98780.5
    "L102" = 98780.5
    "N102" = -105270.5
    "H130" = 92625.14
    "J130" = 92625.14
    "L130" = 92625.14
    "N130" = -102665.14
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H100" = 3008.7144
    "I100" = 3073.6667
    "K100" = 6147.3334
    "M100" = -5606.3334
    "H113" = 1161.8
    "I113" = 471.0625
    "J113" = 3924.75
    "K113" = 1413.1875
    "L113" = 11774.25
    "M113" = 756.8125
    "N113" = -16114.25
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Special: remove M46 cell entirely on LTW ----
$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("M46").ClearContents()

Write-Output "Done applying updates."